$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D36").Value = "Introduction to Human Pose Estimation"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/311"

$ws.Range("D39").Value = "Convolutional Autoencoder: Clustering Images with Neural Networks"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Convolutional-Autoencoder-Clustering-Images-with-Neural-Networks-1"

$ws.Range("D51").Value = "[세이버메트릭스] 타율이 높은 팀 vs OPS가 높은 팀, 누가 이길까?"
$ws.Range("E51").Value = "https://bskyvision.com/1121"

$wb.Save()
